$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '68.354.78'
$c.Style = "Normal"

$c = $ws.Range('E2')
$c.NumberFormat = "@"
$c.Value = '  +1.53%  '
$c.Style = "Normal"

$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '3.899.79'
$c.Style = "Normal"

$c = $ws.Range('E3')
$c.NumberFormat = "@"
$c.Value = '  +1.10%  '
$c.Style = "Normal"

$c = $ws.Range('E4')
$c.NumberFormat = "@"
$c.Value = '  +0.22%  '
$c.Style = "Normal"

$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '484.08'
$c.Style = "Normal"

$c = $ws.Range('E5')
$c.NumberFormat = "@"
$c.Value = '  +4.03%  '
$c.Style = "Normal"

$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '146.66'
$c.Style = "Normal"

$c = $ws.Range('E6')
$c.NumberFormat = "@"
$c.Value = '  -1.45%  '
$c.Style = "Normal"

$c = $ws.Range('E7')
$c.NumberFormat = "@"
$c.Value = '  -2.27%  '
$c.Style = "Normal"

$c = $ws.Range('E8')
$c.NumberFormat = "@"
$c.Value = '  +0.00%  '
$c.Style = "Normal"

$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.718'
$c.Style = "Normal"

$c = $ws.Range('E9')
$c.NumberFormat = "@"
$c.Value = '  -4.27%  '
$c.Style = "Normal"

$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.167'
$c.Style = "Normal"

$c = $ws.Range('E10')
$c.NumberFormat = "@"
$c.Value = '  +8.28%  '
$c.Style = "Normal"

$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.0000352'
$c.Style = "Normal"

$c = $ws.Range('E11')
$c.NumberFormat = "@"
$c.Value = '  +12.95%  '
$c.Style = "Normal"

$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '42.28'
$c.Style = "Normal"

$c = $ws.Range('E12')
$c.NumberFormat = "@"
$c.Value = '  -3.60%  '
$c.Style = "Normal"

$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '10.43'
$c.Style = "Normal"

$c = $ws.Range('E13')
$c.NumberFormat = "@"
$c.Value = '  +0.16%  '
$c.Style = "Normal"

$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '4.524.68'
$c.Style = "Normal"

$c = $ws.Range('E14')
$c.NumberFormat = "@"
$c.Value = '  +1.02%  '
$c.Style = "Normal"

$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '14.51'
$c.Style = "Normal"

$c = $ws.Range('E15')
$c.NumberFormat = "@"
$c.Value = '  -1.44%  '
$c.Style = "Normal"

$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '3.901.17'
$c.Style = "Normal"

$c = $ws.Range('E16')
$c.NumberFormat = "@"
$c.Value = '  +0.20%  '
$c.Style = "Normal"

$c = $ws.Range('E17')
$c.NumberFormat = "@"
$c.Value = '  -0.48%  '
$c.Style = "Normal"

$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '19.65'
$c.Style = "Normal"

$c = $ws.Range('E18')
$c.NumberFormat = "@"
$c.Value = '  -2.02%  '
$c.Style = "Normal"

$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '1.12'
$c.Style = "Normal"

$c = $ws.Range('E19')
$c.NumberFormat = "@"
$c.Value = '  -3.73%  '
$c.Style = "Normal"

$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '68.519.43'
$c.Style = "Normal"

$c = $ws.Range('E20')
$c.NumberFormat = "@"
$c.Value = '  +1.62%  '
$c.Style = "Normal"

$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '430.34'
$c.Style = "Normal"

$c = $ws.Range('E21')
$c.NumberFormat = "@"
$c.Value = '  -0.17%  '
$c.Style = "Normal"

$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '14.47'
$c.Style = "Normal"

$c = $ws.Range('E22')
$c.NumberFormat = "@"
$c.Value = '  -2.37%  '
$c.Style = "Normal"

$c = $ws.Range('E23')
$c.NumberFormat = "@"
$c.Value = '  +1.25%  '
$c.Style = "Normal"

$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '86.63'
$c.Style = "Normal"

$c = $ws.Range('E24')
$c.NumberFormat = "@"
$c.Value = '  -1.98%  '
$c.Style = "Normal"

$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '11.49'
$c.Style = "Normal"

$c = $ws.Range('E25')
$c.NumberFormat = "@"
$c.Value = '  +12.63%  '
$c.Style = "Normal"

$c = $ws.Range('B26')
$c.NumberFormat = "@"
$c.Value = 'PancakeSwap'
$c.Style = "Normal"

$c = $ws.Range('C26')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c.Style = "Normal"

$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '3.57'
$c.Style = "Normal"

$c = $ws.Range('E26')
$c.NumberFormat = "@"
$c.Value = '  +0.46%  '
$c.Style = "Normal"

$c = $ws.Range('B27')
$c.NumberFormat = "@"
$c.Value = 'RenderToken'
$c.Style = "Normal"

$c = $ws.Range('C27')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c.Style = "Normal"

$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '10.57'
$c.Style = "Normal"

$c = $ws.Range('E27')
$c.NumberFormat = "@"
$c.Value = '  +1.42%  '
$c.Style = "Normal"

$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '37.82'
$c.Style = "Normal"

$c = $ws.Range('E28')
$c.NumberFormat = "@"
$c.Value = '  +0.78%  '
$c.Style = "Normal"

$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '5.86'
$c.Style = "Normal"

$c = $ws.Range('E29')
$c.NumberFormat = "@"
$c.Value = '  +6.79%  '
$c.Style = "Normal"

$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '713.09'
$c.Style = "Normal"

$c = $ws.Range('E30')
$c.NumberFormat = "@"
$c.Value = '  -4.36%  '
$c.Style = "Normal"

$c = $ws.Range('E31')
$c.NumberFormat = "@"
$c.Value = '  -4.37%  '
$c.Style = "Normal"

$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '13.17'
$c.Style = "Normal"

$c = $ws.Range('E32')
$c.NumberFormat = "@"
$c.Value = '  -4.26%  '
$c.Style = "Normal"

$c = $ws.Range('E33')
$c.NumberFormat = "@"
$c.Value = '  +2.63%  '
$c.Style = "Normal"

$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '0.0₃0889'
$c.Style = "Normal"

$c = $ws.Range('E34')
$c.NumberFormat = "@"
$c.Value = '  +31.51%  '
$c.Style = "Normal"

$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '41.20'
$c.Style = "Normal"

$c = $ws.Range('E35')
$c.NumberFormat = "@"
$c.Value = '  -4.94%  '
$c.Style = "Normal"

$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '57.86'
$c.Style = "Normal"

$c = $ws.Range('E36')
$c.NumberFormat = "@"
$c.Value = '  +0.93%  '
$c.Style = "Normal"

$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '0.150'
$c.Style = "Normal"

$c = $ws.Range('E37')
$c.NumberFormat = "@"
$c.Value = '  -7.08%  '
$c.Style = "Normal"

$c = $ws.Range('E38')
$c.NumberFormat = "@"
$c.Value = '  -1.65%  '
$c.Style = "Normal"

$c = $ws.Range('E39')
$c.NumberFormat = "@"
$c.Value = '  -0.10%  '
$c.Style = "Normal"

$c = $ws.Range('B40')
$c.NumberFormat = "@"
$c.Value = 'VeChain'
$c.Style = "Normal"

$c = $ws.Range('C40')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c.Style = "Normal"

$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '0.0468'
$c.Style = "Normal"

$c = $ws.Range('E40')
$c.NumberFormat = "@"
$c.Value = '  -2.18%  '
$c.Style = "Normal"

$c = $ws.Range('B41')
$c.NumberFormat = "@"
$c.Value = 'Fetch.AI'
$c.Style = "Normal"

$c = $ws.Range('C41')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c.Style = "Normal"

$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '2.79'
$c.Style = "Normal"

$c = $ws.Range('E41')
$c.NumberFormat = "@"
$c.Value = '  +5.22%  '
$c.Style = "Normal"

$c = $ws.Range('E42')
$c.NumberFormat = "@"
$c.Value = '  +10.92%  '
$c.Style = "Normal"

$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '2.96'
$c.Style = "Normal"

$c = $ws.Range('E43')
$c.NumberFormat = "@"
$c.Value = '  +1.48%  '
$c.Style = "Normal"

$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '0.340'
$c.Style = "Normal"

$c = $ws.Range('E44')
$c.NumberFormat = "@"
$c.Value = '  -3.59%  '
$c.Style = "Normal"

$c = $ws.Range('B45')
$c.NumberFormat = "@"
$c.Value = 'Stellar'
$c.Style = "Normal"

$c = $ws.Range('C45')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c.Style = "Normal"

$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '0.140'
$c.Style = "Normal"

$c = $ws.Range('E45')
$c.NumberFormat = "@"
$c.Value = '  -1.12%  '
$c.Style = "Normal"

$c = $ws.Range('B46')
$c.NumberFormat = "@"
$c.Value = 'FirstDigitalUSD'
$c.Style = "Normal"

$c = $ws.Range('C46')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$c.Style = "Normal"

$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"

$c = $ws.Range('E46')
$c.NumberFormat = "@"
$c.Value = '  -0.09%  '
$c.Style = "Normal"

$c = $ws.Range('B47')
$c.NumberFormat = "@"
$c.Value = 'LidoDAOToken'
$c.Style = "Normal"

$c = $ws.Range('C47')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c.Style = "Normal"

$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '3.39'
$c.Style = "Normal"

$c = $ws.Range('E47')
$c.NumberFormat = "@"
$c.Value = '  -1.56%  '
$c.Style = "Normal"

$c = $ws.Range('B48')
$c.NumberFormat = "@"
$c.Value = 'ARBITRUM'
$c.Style = "Normal"

$c = $ws.Range('C48')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c.Style = "Normal"

$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '2.15'
$c.Style = "Normal"

$c = $ws.Range('E48')
$c.NumberFormat = "@"
$c.Value = '  +0.94%  '
$c.Style = "Normal"

$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '148.22'
$c.Style = "Normal"

$c = $ws.Range('E49')
$c.NumberFormat = "@"
$c.Value = '  +2.86%  '
$c.Style = "Normal"

$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '3.18'
$c.Style = "Normal"

$c = $ws.Range('E50')
$c.NumberFormat = "@"
$c.Value = '  -2.71%  '
$c.Style = "Normal"

$c = $ws.Range('E51')
$c.NumberFormat = "@"
$c.Value = '  -2.91%  '
$c.Style = "Normal"

